$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "a"
$ws.Range("B3").Value = "b1"
$ws.Range("B4").Value = "b2"
$ws.Range("B5").Value = "c1"
$ws.Range("B6").Value = "c2"
$ws.Range("B7").Value = "c3"
$ws.Range("B8").Value = "c4"
